$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Rule Name text for row 3 (paragraph -> wider name)
$ws.Range("B3").Value = "paragraph-with-a-really-wide-rule-name"

# Update the Description text for row 4 so it is wide enough to need wrapping
$ws.Range("D4").Value = "Here's a one line description but it is very wide so should wrap within a cell."

# Wrap text for the whole Rule Name column data range (B2:B6)
$ws.Range("B2:B6").WrapText = $true

# Widen the Rule Name (B) and Description (D) columns, and stop auto bestFit sizing.
# (Target stored widths are 20.7109375 / 30.7109375 "characters"; this runtime's
# column-width engine quantizes to whole pixels, so 19.8 / 29.8 land in the pixel
# bucket closest to those targets.)
$ws.Columns.Item(2).ColumnWidth = 19.8
$ws.Columns.Item(4).ColumnWidth = 29.8
